# temp solve of RWheel
# Set the "Fitness" column (C) values for rows 2-12 to a constant 3917.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2:C12").Value = 3917
